$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "41.569.05"
$ws.Range("E2").Value = "  +0.04%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.210.10"
$ws.Range("E3").Value = "  -1.93%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "228.99"
$ws.Range("E5").Value = "  -1.37%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.616"
$ws.Range("E6").Value = "  -3.28%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "59.50"
$ws.Range("E7").Value = "  -7.33%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.400"
$ws.Range("E9").Value = "  -2.43%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "57.55"
$ws.Range("E10").Value = "  -2.89%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0888"
$ws.Range("E11").Value = "  -1.34%  "
$ws.Range("E12").Value = "  -1.31%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.541.97"
$ws.Range("E13").Value = "  -1.92%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "15.35"
$ws.Range("E14").Value = "  -5.84%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "22.24"
$ws.Range("E15").Value = "  -1.52%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.62"
$ws.Range("E16").Value = "  -1.29%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.790"
$ws.Range("E17").Value = "  -5.17%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.228.24"
$ws.Range("E18").Value = "  -0.81%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "41.487.32"
$ws.Range("E19").Value = "  +0.14%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0897"
$ws.Range("E20").Value = "  -1.87%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "71.68"
$ws.Range("E21").Value = "  -3.14%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.05"
$ws.Range("E22").Value = "  -2.55%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "241.84"
$ws.Range("E23").Value = "  -4.00%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.998"
$ws.Range("E24").Value = "  -0.26%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.34"
$ws.Range("E25").Value = "  -2.40%  "
$ws.Range("E26").Value = "  -5.23%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.65"
$ws.Range("E27").Value = "  -2.03%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "168.44"
$ws.Range("E28").Value = "  -2.96%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.139"
$ws.Range("E29").Value = "  -4.64%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "19.69"
$ws.Range("E30").Value = "  -3.91%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.40"
$ws.Range("E31").Value = "  -4.32%  "
$ws.Range("E32").Value = "  -10.20%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.121"
$ws.Range("E33").Value = "  -3.32%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.93"
$ws.Range("E34").Value = "  -2.06%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.62"
$ws.Range("E35").Value = "  -2.73%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0647"
$ws.Range("E36").Value = "  +1.62%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.44"
$ws.Range("E37").Value = "  -8.50%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.35"
$ws.Range("E38").Value = "  -4.49%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.55"
$ws.Range("E39").Value = "  -7.50%  "
$ws.Range("E40").Value = "  +0.04%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.000231"
$ws.Range("E41").Value = "  -13.56%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0236"
$ws.Range("E42").Value = "  -1.70%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.55"
$ws.Range("E43").Value = "  -3.41%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0958"
$ws.Range("E44").Value = "  +1.42%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.20"
$ws.Range("E45").Value = "  -3.06%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "96.90"
$ws.Range("E46").Value = "  -5.73%  "
$ws.Range("B47").Value = "Maker"
$ws.Range("C47").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.463.07"
$ws.Range("E47").Value = "  -3.17%  "
$ws.Range("B48").Value = "FTXToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "4.35"
$ws.Range("E48").Value = "  -10.75%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "16.36"
$ws.Range("E49").Value = "  -8.47%  "
$ws.Range("E50").Value = "  -1.59%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.06"
$ws.Range("E51").Value = "  -5.31%  "
